$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.933.00"
$ws.Range("D3").Value = "3.321.66"
$ws.Range("E3").Value = "  +6.31%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.16%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.320.83"
$ws.Range("E8").Value = "  +6.80%  "
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("E10").Value = "  +3.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.474"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.60%  "
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").Value = "3.866.18"
$ws.Range("E15").Value = "  +6.30%  "
$ws.Range("D17").Value = "3.320.90"
$ws.Range("E17").Value = "  +6.47%  "
$ws.Range("D18").Value = "64.013.18"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.738"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.44%  "
$ws.Range("E23").Value = "  +3.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.84%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.80%  "
$ws.Range("E31").Value = "  +5.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.98%  "
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.20%  "
$ws.Range("D37").Value = "0.0₃0758"
$ws.Range("E37").Value = "  +8.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0404"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "432.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("D41").Value = "3.063.26"
$ws.Range("E41").Value = "  +5.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("E46").Value = "  +5.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.115"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.62%  "
$ws.Range("E51").Value = "  +2.90%  "

Write-Output "Applied 76 changes"